$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that are stored as plain text in the source data
# (e.g. "34.503.97", "0.603"). Assigning such numeric-looking strings directly to
# .Value would make Excel auto-convert them to numbers, so a leading apostrophe is
# used to force a text value (Excel strips the apostrophe from the stored value).

# Row 2
$ws.Range("D2").Value = "34.503.97"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").Value = "1.801.50"
$ws.Range("E3").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'224.59"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").Value = "'0.603"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "'42.37"
$ws.Range("E8").Value = "  +17.82%  "

# Row 9
$ws.Range("D9").Value = "'0.292"
$ws.Range("E9").Value = "  +0.38%  "

# Row 10
$ws.Range("D10").Value = "'0.0666"
$ws.Range("E10").Value = "  -1.64%  "

# Row 11
$ws.Range("E11").Value = "  +3.10%  "

# Row 12
$ws.Range("D12").Value = "2.062.08"
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$ws.Range("D13").Value = "1.800.85"
$ws.Range("E13").Value = "  -0.10%  "

# Row 14
$ws.Range("D14").Value = "'10.93"
$ws.Range("E14").Value = "  -2.49%  "

# Row 15
$ws.Range("D15").Value = "34.463.50"
$ws.Range("E15").Value = "  +0.28%  "

# Row 16
$ws.Range("D16").Value = "'0.628"
$ws.Range("E16").Value = "  -0.27%  "

# Row 17
$ws.Range("D17").Value = "'4.41"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
$ws.Range("E18").Value = "  -1.81%  "

# Row 19
$ws.Range("D19").Value = "'240.27"
$ws.Range("E19").Value = "  -0.80%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0768"
$ws.Range("E20").Value = "  -0.65%  "

# Row 21
$ws.Range("D21").Value = "'11.15"
$ws.Range("E21").Value = "  -1.16%  "

# Row 22
$ws.Range("E22").Value = "  +0.18%  "

# Row 23
$ws.Range("D23").Value = "'4.37"
$ws.Range("E23").Value = "  +7.12%  "

# Row 24
$ws.Range("E24").Value = "  -2.77%  "

# Row 25
$ws.Range("D25").Value = "'170.97"
$ws.Range("E25").Value = "  +0.47%  "

# Row 26
$ws.Range("D26").Value = "'7.66"
$ws.Range("E26").Value = "  -2.54%  "

# Row 27
$ws.Range("D27").Value = "'17.40"

# Row 28
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("E31").Value = "  -0.22%  "

# Row 32
$ws.Range("E32").Value = "  -1.11%  "

# Row 33
$ws.Range("D33").Value = "'0.0512"
$ws.Range("E33").Value = "  -0.62%  "

# Row 34
$ws.Range("E34").Value = "  +0.71%  "

# Row 35
$ws.Range("D35").Value = "'87.67"
$ws.Range("E35").Value = "  +7.95%  "

# Row 36
$ws.Range("D36").Value = "'0.647"
$ws.Range("E36").Value = "  -0.46%  "

# Row 37
$ws.Range("D37").Value = "1.318.64"
$ws.Range("E37").Value = "  -3.33%  "

# Row 38
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = "  +0.17%  "

# Row 39
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'14.82"
$ws.Range("E39").Value = "  +12.02%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0188"
$ws.Range("E40").Value = "  +0.86%  "

# Row 41
$ws.Range("E41").Value = "  -1.34%  "

# Row 42
$ws.Range("E42").Value = "  +4.88%  "

# Row 43
$ws.Range("E43").Value = "  +0.26%  "

# Row 44
$ws.Range("E44").Value = "  +0.24%  "

# Row 45
$ws.Range("D45").Value = "'0.936"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("E46").Value = "  +3.66%  "

# Row 47
$ws.Range("D47").Value = "1.963.95"
$ws.Range("E47").Value = "  +0.10%  "

# Row 48
$ws.Range("E48").Value = "  +0.11%  "

# Row 49
$ws.Range("E49").Value = "  +0.14%  "

# Row 50
$ws.Range("D50").Value = "'100.50"
$ws.Range("E50").Value = "  -1.36%  "

# Row 51
$ws.Range("D51").Value = "'0.0609"
$ws.Range("E51").Value = "  +0.77%  "
